# Applies scheduled-runner market-price refresh to the Leve profit tables.
# For each sheet, update the cached currentAveragePrice* / LevePrice* / LeveProfit*
# columns (H..N) with refreshed values; some rows gain or lose a LeveProfitHQ (N) cell
# entirely depending on whether HQ data is available.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 707.7174
$ws.Range("J17").Value = 707.7174
$ws.Range("L17").Value = 2123.1522
$ws.Range("N17").Value = -2459.1522
$ws.Range("H32").Value = 3499.5
$ws.Range("J32").Value = 3499.5
$ws.Range("L32").Value = 3499.5
$ws.Range("N32").Value = -4151.5
$ws.Range("H51").Value = 10980.091
$ws.Range("J51").Value = 9078.1
$ws.Range("L51").Value = 9078.1
$ws.Range("N51").Value = -10046.1
$ws.Range("H86").Value = 3545.077
$ws.Range("I86").Value = 3260.625
$ws.Range("K86").Value = 3260.625
$ws.Range("M86").Value = -2137.625
$ws.Range("H89").Value = 3545.077
$ws.Range("I89").Value = 3260.625
$ws.Range("K89").Value = 16303.125
$ws.Range("M89").Value = -10687.125
$ws.Range("H111").Value = 4661.1
$ws.Range("I111").Value = 4734.5557
$ws.Range("J111").Value = 4000
$ws.Range("K111").Value = 14203.6671
$ws.Range("L111").Value = 12000
$ws.Range("M111").Value = -11136.6671
$ws.Range("N111").Value = -18134
$ws.Range("H113").Value = 6411.579
$ws.Range("I113").Value = 6109.2856
$ws.Range("J113").Value = 6587.9165
$ws.Range("K113").Value = 6109.2856
$ws.Range("L113").Value = 6587.9165
$ws.Range("M113").Value = -2855.2856
$ws.Range("N113").Value = -13095.9165
$ws.Range("H116").Value = 11039.315
$ws.Range("J116").Value = 8943.647000000001
$ws.Range("L116").Value = 8943.647000000001
$ws.Range("N116").Value = -15827.647
$ws.Range("H132").Value = 3396.5557
$ws.Range("I132").Value = 3071.125
$ws.Range("K132").Value = 9213.375
$ws.Range("M132").Value = -6683.375
$ws.Range("H137").Value = 403180.9
$ws.Range("I137").Value = 437894.1
$ws.Range("K137").Value = 1313682.3
$ws.Range("M137").Value = -1311132.3
$ws.Range("H141").Value = 6716.7334
$ws.Range("I141").Value = 4269
$ws.Range("J141").Value = 8858.5
$ws.Range("K141").Value = 12807
$ws.Range("L141").Value = 26575.5
$ws.Range("M141").Value = -7627
$ws.Range("N141").Value = -36935.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1178.138
$ws.Range("I2").Value = 1116.625
$ws.Range("J2").Value = 1473.4
$ws.Range("K2").Value = 1116.625
$ws.Range("L2").Value = 1473.4
$ws.Range("M2").Value = -1003.625
$ws.Range("N2").Value = -1699.4
$ws.Range("H5").Value = 637.9375
$ws.Range("I5").Value = 655.2727
$ws.Range("J5").Value = 599.8
$ws.Range("K5").Value = 655.2727
$ws.Range("L5").Value = 599.8
$ws.Range("M5").Value = -543.2727
$ws.Range("N5").Value = -823.8
$ws.Range("H16").Value = 16499
$ws.Range("I16").Value = 999
$ws.Range("J16").Value = 21665.666
$ws.Range("K16").Value = 999
$ws.Range("L16").Value = 21665.666
$ws.Range("M16").Value = -712
$ws.Range("N16").Value = -22239.666
$ws.Range("H42").Value = 26000
$ws.Range("J42").Value = 26000
$ws.Range("L42").Value = 26000
$ws.Range("N42").Value = -26972
$ws.Range("H44").Value = 52499
$ws.Range("J44").Value = 52499
$ws.Range("L44").Value = 52499
$ws.Range("N44").Value = -53475
$ws.Range("H61").Value = 7815.25
$ws.Range("I61").Value = 3561.9412
$ws.Range("K61").Value = 3561.9412
$ws.Range("M61").Value = -3349.9412
$ws.Range("H80").Value = 60000
$ws.Range("J80").Value = 60000
$ws.Range("L80").Value = 60000
$ws.Range("N80").Value = -61996
$ws.Range("H83").Value = 60000
$ws.Range("J83").Value = 60000
$ws.Range("L83").Value = 180000
$ws.Range("N83").Value = -189984
$ws.Range("H102").Value = 1291.3846
$ws.Range("I102").Value = 1291.3846
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 1291.3846
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = 330.6153999999999
$ws.Range("N102").ClearContents()
$ws.Range("H116").Value = 1178.138
$ws.Range("I116").Value = 1116.625
$ws.Range("J116").Value = 1473.4
$ws.Range("K116").Value = 1116.625
$ws.Range("L116").Value = 1473.4
$ws.Range("M116").Value = 1177.375
$ws.Range("N116").Value = -6061.4
$ws.Range("H132").Value = 1773.625
$ws.Range("I132").Value = 1500.3103
$ws.Range("K132").Value = 4500.9309
$ws.Range("M132").Value = -1970.9309
$ws.Range("H136").Value = 7815.25
$ws.Range("I136").Value = 3561.9412
$ws.Range("K136").Value = 10685.8236
$ws.Range("M136").Value = -8135.8236

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1178.138
$ws.Range("I3").Value = 1116.625
$ws.Range("J3").Value = 1473.4
$ws.Range("K3").Value = 1116.625
$ws.Range("L3").Value = 1473.4
$ws.Range("M3").Value = -1002.625
$ws.Range("N3").Value = -1701.4
$ws.Range("H4").Value = 637.9375
$ws.Range("I4").Value = 655.2727
$ws.Range("J4").Value = 599.8
$ws.Range("K4").Value = 655.2727
$ws.Range("L4").Value = 599.8
$ws.Range("M4").Value = -540.2727
$ws.Range("N4").Value = -829.8
$ws.Range("H82").Value = 29735.5
$ws.Range("J82").Value = 48996
$ws.Range("L82").Value = 48996
$ws.Range("N82").Value = -49762
$ws.Range("H85").Value = 29735.5
$ws.Range("J85").Value = 48996
$ws.Range("L85").Value = 48996
$ws.Range("N85").Value = -51648
$ws.Range("H99").Value = 4561
$ws.Range("J99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("N99").ClearContents()
$ws.Range("H107").Value = 1906.3125
$ws.Range("I107").Value = 2366
$ws.Range("K107").Value = 2366
$ws.Range("M107").Value = -446

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 985.8570999999999
$ws.Range("J22").Value = 825
$ws.Range("L22").Value = 825
$ws.Range("N22").Value = -1525
$ws.Range("H105").Value = 1665.7368
$ws.Range("I105").Value = 1348.1666
$ws.Range("K105").Value = 1348.1666
$ws.Range("M105").Value = 398.8334
$ws.Range("H132").Value = 3381.4666
$ws.Range("I132").Value = 3380.1428
$ws.Range("K132").Value = 10140.4284
$ws.Range("M132").Value = -7610.428400000001
$ws.Range("H134").Value = 5575.8965
$ws.Range("I134").Value = 4218.4546
$ws.Range("J134").Value = 9842.143
$ws.Range("K134").Value = 12655.3638
$ws.Range("L134").Value = 29526.429
$ws.Range("M134").Value = -10120.3638
$ws.Range("N134").Value = -34596.429

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 1860
$ws.Range("H37").Value = 139814.38
$ws.Range("J37").Value = 139814.38
$ws.Range("L37").Value = 419443.14
$ws.Range("N37").Value = -419667.14
$ws.Range("H97").Value = 1232.3334
$ws.Range("I97").Value = 1198
$ws.Range("J97").Value = 1239.2
$ws.Range("K97").Value = 3594
$ws.Range("L97").Value = 3717.6
$ws.Range("M97").Value = -3098
$ws.Range("N97").Value = -4709.6
$ws.Range("H107").Value = 1239.8125
$ws.Range("J107").Value = 1324.8572
$ws.Range("L107").Value = 3974.5716
$ws.Range("N107").Value = -7814.571599999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3113.5
$ws.Range("I80").Value = 3113.5
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 3113.5
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -2115.5
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 3113.5
$ws.Range("I83").Value = 3113.5
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 15567.5
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -10575.5
$ws.Range("N83").ClearContents()
$ws.Range("H102").Value = 1500
$ws.Range("I102").Value = 1500
$ws.Range("K102").Value = 1500
$ws.Range("M102").Value = 122
$ws.Range("H113").Value = 1536
$ws.Range("I113").Value = 964.7
$ws.Range("K113").Value = 964.7
$ws.Range("M113").Value = 1205.3
$ws.Range("H132").Value = 4228.7617
$ws.Range("I132").Value = 4315.2
$ws.Range("K132").Value = 12945.6
$ws.Range("M132").Value = -10415.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2553
$ws.Range("I40").Value = 2230
$ws.Range("K40").Value = 2230
$ws.Range("M40").Value = -2094
$ws.Range("H101").Value = 6746.75
$ws.Range("J101").Value = 6746.75
$ws.Range("L101").Value = 6746.75
$ws.Range("N101").Value = -13236.75
$ws.Range("H136").Value = 2101.925
$ws.Range("I136").Value = 1257.1666
$ws.Range("J136").Value = 3369.0625
$ws.Range("K136").Value = 3771.4998
$ws.Range("L136").Value = 10107.1875
$ws.Range("M136").Value = -1221.4998
$ws.Range("N136").Value = -15207.1875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 12500
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()
$ws.Range("H113").Value = 1486
$ws.Range("I113").Value = 1814.6666
$ws.Range("J113").Value = 500
$ws.Range("K113").Value = 5443.9998
$ws.Range("L113").Value = 1500
$ws.Range("M113").Value = -3273.9998
$ws.Range("N113").Value = -5840
$ws.Range("H122").Value = 5414.2856
$ws.Range("I122").Value = 5414.2856
$ws.Range("K122").Value = 16242.8568
$ws.Range("M122").Value = -13792.8568
$ws.Range("H136").Value = 5706.653
$ws.Range("I136").Value = 5067.978
$ws.Range("J136").Value = 15499.667
$ws.Range("K136").Value = 15203.934
$ws.Range("L136").Value = 46499.001
$ws.Range("M136").Value = -12653.934
$ws.Range("N136").Value = -51599.001
